# Commit: "added user manager y recuperacion de usuarrio"
#
# Changes:
#  1. Inventario (sheet1): Fernet Branca stock corrected 108 -> 98, and a new
#     product row (id 8, "nalga de tom") is added as row 9.
#  2. Historia (sheet2): five new movement/log rows are appended (rows 7-11):
#     - two stock egresos (Fernet Branca, nalga de tom)
#     - three "Inicio de sesión" (login) events with placeholder "-" product
#       and zeroed quantity/price/total, reflecting the new user-manager /
#       password-recovery + login-logging feature.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Inventario")
$ws2 = $wb.Worksheets.Item("Historia")

# --- Inventario: fix existing stock value -------------------------------
$ws1.Range("C3").Value = 98

# --- Inventario: add new product row (row 9) -----------------------------
$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "nalga de tom"
$ws1.Range("C9").Value = 99984
$ws1.Range("D9").Value = 115

# --- Historia: append new movement rows -----------------------------------
$ws2.Range("A7").Value = "2025-08-04 15:31:33"
$ws2.Range("B7").Value = "blestro"
$ws2.Range("C7").Value = "Fernet Branca"
$ws2.Range("D7").Value = -10
$ws2.Range("E7").Value = 4500
$ws2.Range("F7").Value = 45000
$ws2.Range("G7").Value = "Egreso"

$ws2.Range("A8").Value = "2025-08-04 15:32:06"
$ws2.Range("B8").Value = "blestro"
$ws2.Range("C8").Value = "nalga de tom"
$ws2.Range("D8").Value = -16
$ws2.Range("E8").Value = 115
$ws2.Range("F8").Value = 1840
$ws2.Range("G8").Value = "Egreso"

$ws2.Range("A9").Value = "2025-08-04 18:24:44"
$ws2.Range("B9").Value = "blestro"
$ws2.Range("C9").Value = "-"
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = "Inicio de sesión"

$ws2.Range("A10").Value = "2025-08-04 18:40:02"
$ws2.Range("B10").Value = "blestro"
$ws2.Range("C10").Value = "-"
$ws2.Range("D10").Value = 0
$ws2.Range("E10").Value = 0
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = "Inicio de sesión"

$ws2.Range("A11").Value = "2025-08-04 18:46:00"
$ws2.Range("B11").Value = "blestro"
$ws2.Range("C11").Value = "-"
$ws2.Range("D11").Value = 0
$ws2.Range("E11").Value = 0
$ws2.Range("F11").Value = 0
$ws2.Range("G11").Value = "Inicio de sesión"
